# Generate Report for Handoff
# - Mark the 6 "ready for handoff" rows (4bb9a067, 63e9217c, 65ff044d,
#   8142fc67, b2a25f4e, bb13e498) with Priority "ht" in both the zh-cn and
#   de-de localization sheets.
# - Refresh the "Latest Handoff"/"Latest HO Xliff Generate Date" timestamps
#   for those same rows across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 11, 12, 13, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    # Overview sheet: column G = "Latest HO Xliff Generate Date"
    $wsOverview.Range("G$r").Value = "2016-08-21 04:20:20"

    # zh-cn sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-21 04:20:15"

    # de-de sheet: column E = "Priority", column H = "Latest Handoff Datetime"
    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-21 04:20:20"
}
